# Refresh the crypto price/volume table to the latest scraped snapshot.
# (GitHub Actions commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.414.80'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '1.639.80'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').Formula = '="304.67"'
$ws.Range('D6').Copy() | Out-Null
$ws.Range('D6').PasteSpecial(-4163) | Out-Null
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('D7').Formula = '="0.3726"'
$ws.Range('D7').Copy() | Out-Null
$ws.Range('D7').PasteSpecial(-4163) | Out-Null
$ws.Range('E7').Value = '  -1.11%  '
$ws.Range('D8').Formula = '="52.40"'
$ws.Range('D8').Copy() | Out-Null
$ws.Range('D8').PasteSpecial(-4163) | Out-Null
$ws.Range('E8').Value = '  +1.10%  '
$ws.Range('D9').Formula = '="0.3623"'
$ws.Range('D9').Copy() | Out-Null
$ws.Range('D9').PasteSpecial(-4163) | Out-Null
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('D10').Formula = '="1.248"'
$ws.Range('D10').Copy() | Out-Null
$ws.Range('D10').PasteSpecial(-4163) | Out-Null
$ws.Range('E10').Value = '  -2.16%  '
$ws.Range('D11').Formula = '="0.08103"'
$ws.Range('D11').Copy() | Out-Null
$ws.Range('D11').PasteSpecial(-4163) | Out-Null
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Formula = '="22.80"'
$ws.Range('D13').Copy() | Out-Null
$ws.Range('D13').PasteSpecial(-4163) | Out-Null
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('D14').Formula = '="6.587"'
$ws.Range('D14').Copy() | Out-Null
$ws.Range('D14').PasteSpecial(-4163) | Out-Null
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').Formula = '="0.00001269"'
$ws.Range('D15').Copy() | Out-Null
$ws.Range('D15').PasteSpecial(-4163) | Out-Null
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('D16').Formula = '="7.281"'
$ws.Range('D16').Copy() | Out-Null
$ws.Range('D16').PasteSpecial(-4163) | Out-Null
$ws.Range('E16').Value = '  -1.89%  '
$ws.Range('D17').Value = '1.631.25'
$ws.Range('E17').Value = '  +1.54%  '
$ws.Range('D18').Formula = '="94.40"'
$ws.Range('D18').Copy() | Out-Null
$ws.Range('D18').PasteSpecial(-4163) | Out-Null
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('D19').Formula = '="0.06888"'
$ws.Range('D19').Copy() | Out-Null
$ws.Range('D19').PasteSpecial(-4163) | Out-Null
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').Formula = '="6.504"'
$ws.Range('D21').Copy() | Out-Null
$ws.Range('D21').PasteSpecial(-4163) | Out-Null
$ws.Range('E21').Value = '  -0.37%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '23.427.99'
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('D24').Formula = '="12.74"'
$ws.Range('D24').Copy() | Out-Null
$ws.Range('D24').PasteSpecial(-4163) | Out-Null
$ws.Range('E24').Value = '  -1.63%  '
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').Formula = '="3.059"'
$ws.Range('D25').Copy() | Out-Null
$ws.Range('D25').PasteSpecial(-4163) | Out-Null
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Formula = '="2.411"'
$ws.Range('D26').Copy() | Out-Null
$ws.Range('D26').PasteSpecial(-4163) | Out-Null
$ws.Range('E26').Value = '  +1.29%  '
$ws.Range('D27').Formula = '="21.17"'
$ws.Range('D27').Copy() | Out-Null
$ws.Range('D27').PasteSpecial(-4163) | Out-Null
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').Formula = '="151.36"'
$ws.Range('D28').Copy() | Out-Null
$ws.Range('D28').PasteSpecial(-4163) | Out-Null
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('D29').Formula = '="5.329"'
$ws.Range('D29').Copy() | Out-Null
$ws.Range('D29').PasteSpecial(-4163) | Out-Null
$ws.Range('E29').Value = '  +1.42%  '
$ws.Range('D30').Formula = '="135.71"'
$ws.Range('D30').Copy() | Out-Null
$ws.Range('D30').PasteSpecial(-4163) | Out-Null
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('D31').Formula = '="2.282"'
$ws.Range('D31').Copy() | Out-Null
$ws.Range('D31').PasteSpecial(-4163) | Out-Null
$ws.Range('E31').Value = '  -4.90%  '
$ws.Range('D32').Value = '1.810.86'
$ws.Range('E32').Value = '  +1.51%  '
$ws.Range('D33').Formula = '="6.785"'
$ws.Range('D33').Copy() | Out-Null
$ws.Range('D33').PasteSpecial(-4163) | Out-Null
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('D34').Formula = '="0.9494"'
$ws.Range('D34').Copy() | Out-Null
$ws.Range('D34').PasteSpecial(-4163) | Out-Null
$ws.Range('E34').Value = '  -1.24%  '
$ws.Range('D35').Formula = '="0.02812"'
$ws.Range('D35').Copy() | Out-Null
$ws.Range('D35').PasteSpecial(-4163) | Out-Null
$ws.Range('E35').Value = '  +2.39%  '
$ws.Range('D36').Formula = '="10.35"'
$ws.Range('D36').Copy() | Out-Null
$ws.Range('D36').PasteSpecial(-4163) | Out-Null
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('D38').Formula = '="0.07206"'
$ws.Range('D38').Copy() | Out-Null
$ws.Range('D38').PasteSpecial(-4163) | Out-Null
$ws.Range('E38').Value = '  -3.93%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Formula = '="6.109"'
$ws.Range('D39').Copy() | Out-Null
$ws.Range('D39').PasteSpecial(-4163) | Out-Null
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Formula = '="0.08747"'
$ws.Range('D40').Copy() | Out-Null
$ws.Range('D40').PasteSpecial(-4163) | Out-Null
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').Formula = '="1.370"'
$ws.Range('D41').Copy() | Out-Null
$ws.Range('D41').PasteSpecial(-4163) | Out-Null
$ws.Range('E41').Value = '  -2.02%  '
$ws.Range('D42').Formula = '="0.7022"'
$ws.Range('D42').Copy() | Out-Null
$ws.Range('D42').PasteSpecial(-4163) | Out-Null
$ws.Range('E42').Value = '  -1.27%  '
$ws.Range('E43').Value = '  -0.60%  '
$ws.Range('D44').Formula = '="16.02"'
$ws.Range('D44').Copy() | Out-Null
$ws.Range('D44').PasteSpecial(-4163) | Out-Null
$ws.Range('E44').Value = '  +2.61%  '
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('D46').Formula = '="2.326"'
$ws.Range('D46').Copy() | Out-Null
$ws.Range('D46').PasteSpecial(-4163) | Out-Null
$ws.Range('E46').Value = '  +0.36%  '
$ws.Range('D47').Formula = '="0.9994"'
$ws.Range('D47').Copy() | Out-Null
$ws.Range('D47').PasteSpecial(-4163) | Out-Null
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('D48').Formula = '="4.005"'
$ws.Range('D48').Copy() | Out-Null
$ws.Range('D48').PasteSpecial(-4163) | Out-Null
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('D49').Formula = '="0.07964"'
$ws.Range('D49').Copy() | Out-Null
$ws.Range('D49').PasteSpecial(-4163) | Out-Null
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').Formula = '="128.28"'
$ws.Range('D50').Copy() | Out-Null
$ws.Range('D50').PasteSpecial(-4163) | Out-Null
$ws.Range('E50').Value = '  -3.30%  '
$ws.Range('D51').Formula = '="1.196"'
$ws.Range('D51').Copy() | Out-Null
$ws.Range('D51').PasteSpecial(-4163) | Out-Null
$ws.Range('E51').Value = '  -0.81%  '

$excel.CutCopyMode = 0
